$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tâches")

# ---------------------------------------------------------------------------
# 1. Add the three new logged tasks (rows 50-52) on 2021-05-21 (serial 44337)
# ---------------------------------------------------------------------------

# Carry the number formats (date / time) used by the existing data rows down
# onto the new rows so the new cells reuse the same styles instead of minting
# new ones.
$ws.Range("B49:E49").Copy() | Out-Null
$ws.Range("B50:E52").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$rows = @(
    @{ Row = 50; Date = 44337; Start = 0.33333333333333331; End = 0.39930555555555558; Type = "Réalisation"; Task = "Suppression des câbles et amélioration de la visibilité" },
    @{ Row = 51; Date = 44337; Start = 0.40972222222222227; End = 0.4597222222222222;  Type = "Réalisation"; Task = "Inputs et outputs du canevas" },
    @{ Row = 52; Date = 44337; Start = 0.4604166666666667;  End = 0.47638888888888892; Type = "Réalisation"; Task = "synérgie canevas et inputs" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Formula = "=ROW()-1"
    $ws.Cells.Item($row, 2).Value = $r.Date
    $ws.Cells.Item($row, 3).Value = $r.Start
    $ws.Cells.Item($row, 4).Value = $r.End
    $ws.Cells.Item($row, 5).Formula = "=IF(OR(ISBLANK(D$row),ISBLANK(C$row)),""."",D$row-C$row)"
    $ws.Cells.Item($row, 6).Value = $r.Type
    $ws.Cells.Item($row, 7).Value = $r.Task
}

# ---------------------------------------------------------------------------
# 2. Extend the Index helper column (rows 53-62) - no other data on those rows
# ---------------------------------------------------------------------------
for ($row = 53; $row -le 62; $row++) {
    $ws.Cells.Item($row, 1).Formula = "=ROW()-1"
}

# ---------------------------------------------------------------------------
# 3. Refresh the pivot tables / pivot caches so they pick up the new records
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $pivotCount = $sheet.PivotTables().Count
    for ($i = 1; $i -le $pivotCount; $i++) {
        $sheet.PivotTables($i).RefreshTable() | Out-Null
    }
}

# ---------------------------------------------------------------------------
# 4. Update the view state on the "Tâches" sheet (scrolled down to row 34,
#    selection on G58) to match where the author ended up after the edit.
# ---------------------------------------------------------------------------
$ws.Activate()
$appWin = $excel.ActiveWindow
$appWin.ScrollRow = 34
$ws.Range("G58").Select()
